$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Login with valid username and password", "PASSED", "chrome", "12_06_23_191833"),
    @("Create Country", "PASSED", "chrome", "12_06_23_191842"),
    @("Add citizenship", "PASSED", "chrome", "12_06_23_191852"),
    @("Add citizenship", "PASSED", "chrome", "12_06_23_191903"),
    @("Add citizenship", "PASSED", "chrome", "12_06_23_191914"),
    @("Add citizenship", "PASSED", "chrome", "12_06_23_191925"),
    @("Login with valid username and password", "PASSED", "chrome", "22_06_23_043238"),
    @("Create Country", "PASSED", "chrome", "22_06_23_043247")
)

$startRow = 58
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}
